$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.258.01"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").Value = "1.896.76"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "'246.49"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").Value = "'0.691"
$ws.Range("E6").Value = "  +9.27%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'40.50"
$ws.Range("E8").Value = "  -4.02%  "

$ws.Range("E9").Value = "  +2.47%  "

$ws.Range("D10").Value = "'52.35"
$ws.Range("E10").Value = "  +8.62%  "

$ws.Range("D11").Value = "'0.0720"
$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("D12").Value = "'0.0985"
$ws.Range("E12").Value = "  -1.30%  "

$ws.Range("D13").Value = "2.171.44"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "'12.52"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D16").Value = "1.915.40"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "'4.80"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").Value = "35.243.71"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").Value = "'71.97"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "0.0₃0819"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").Value = "'240.69"
$ws.Range("E21").Value = "  -1.33%  "

$ws.Range("D22").Value = "'12.70"
$ws.Range("E22").Value = "  +1.67%  "

$ws.Range("D23").Value = "'4.79"
$ws.Range("E23").Value = "  -2.08%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").Value = "'2.32"
$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("D26").Value = "'2.36"
$ws.Range("E26").Value = "  +6.89%  "

$ws.Range("D27").Value = "'168.33"
$ws.Range("E27").Value = "  -2.22%  "

$ws.Range("D28").Value = "'8.64"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").Value = "'18.77"
$ws.Range("E29").Value = "  +4.41%  "

$ws.Range("D30").Value = "'0.131"
$ws.Range("E30").Value = "  +4.57%  "

$ws.Range("E31").Value = "  +20.04%  "

$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("D33").Value = "'0.0567"
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").Value = "'1.86"
$ws.Range("E35").Value = "  +6.56%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.52"
$ws.Range("E36").Value = "  +14.90%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'4.10"
$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("E38").Value = "  -8.45%  "

$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("D40").Value = "'0.0656"
$ws.Range("E40").Value = "  +10.66%  "

$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'16.36"
$ws.Range("E42").Value = "  +5.60%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0207"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").Value = "'92.65"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("D45").Value = "1.348.77"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("E46").Value = "  +2.69%  "

$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").Value = "'2.79"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("B49").Value = "Gas"
$ws.Range("C49").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D49").Value = "'12.49"
$ws.Range("E49").Value = "  -0.87%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'45.32"
$ws.Range("E50").Value = "  -11.14%  "

$ws.Range("E51").Value = "  -3.25%  "
